# Weekly refresh of "Fruta / hortaliza" data: updates Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Unidad de comercializacion, Precio $/Kg and
# Kg o Unidades for every data row (rows 2-28) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is: RowNumber, Fecha(D), Calidad(I), Volumen(J),
# PrecioMinimo(K), PrecioMaximo(L), PrecioPromedioPonderado(M),
# UnidadComercializacion(N), PrecioKg(P), KgOUnidades(Q)
$rowsData = @(
        @(2, 45028, 'Primera', 300, 14000, 15000, 14500, '$/caja 13 kilos', 1115, 13),
        @(3, 44943, 'Segunda', 350, 14000, 15000, 14429, '$/caja 13 kilos', 1110, 13),
        @(4, 44616, 'Primera', 120, 19000, 20000, 19500, '$/caja 13 kilos', 1500, 13),
        @(5, 44580, 'Primera', 160, 11000, 12000, 11500, '$/caja 13 kilos', 885, 13),
        @(6, 44984, 'Primera', 400, 16000, 17000, 16500, '$/caja 13 kilos', 1269, 13),
        @(7, 45100, 'Primera', 200, 15000, 16000, 15500, '$/caja 13 kilos', 1192, 13),
        @(8, 44988, 'Primera', 750, 17000, 18000, 17400, '$/caja 13 kilos', 1338, 13),
        @(9, 44893, 'Primera', 900, 13000, 14000, 13444, '$/caja 13 kilos', 1034, 13),
        @(10, 44592, 'Primera', 120, 12000, 13000, 12500, '$/caja 13 kilos', 962, 13),
        @(11, 44406, 'Primera', 160, 17000, 18000, 17500, '$/caja 13 kilos', 1346, 13),
        @(12, 44389, 'Primera', 120, 12000, 13000, 12500, '$/caja 13 kilos', 962, 13),
        @(13, 44469, 'Primera', 140, 13000, 14000, 13500, '$/caja 13 kilos', 1038, 13),
        @(14, 44764, 'Primera', 200, 12000, 13000, 12500, '$/caja 13 kilos', 962, 13),
        @(15, 44914, 'Primera', 100, 14000, 15000, 14400, '$/caja 13 kilos', 1108, 13),
        @(16, 44159, 'Primera', 100, 23000, 24000, 23500, '$/caja 13 kilos', 1808, 13),
        @(17, 44379, 'Primera', 120, 12000, 13000, 12667, '$/caja 13 kilos', 974, 13),
        @(18, 45092, 'Primera', 600, 13000, 14000, 13500, '$/caja 13 kilos', 1038, 13),
        @(19, 44832, 'Primera', 100, 13000, 14000, 13500, '$/caja 13 kilos', 1038, 13),
        @(20, 44320, 'Primera', 160, 19000, 20000, 19500, '$/caja 13 kilos', 1500, 13),
        @(21, 44910, 'Primera', 50, 14000, 15000, 14500, '$/caja 13 kilos', 1115, 13),
        @(22, 44397, 'Primera', 140, 12500, 13000, 12750, '$/caja 13 kilos', 981, 13),
        @(23, 44890, 'Primera', 300, 14000, 15000, 14500, '$/caja 13 kilos', 1115, 13),
        @(24, 44972, 'Primera', 350, 17000, 18000, 17429, '$/caja 15 kilos', 1162, 15),
        @(25, 44855, 'Primera', 500, 10000, 10000, 10000, '$/caja 13 kilos', 769, 13),
        @(26, 44918, 'Segunda', 200, 12000, 13000, 12750, '$/caja 13 kilos', 981, 13),
        @(27, 45096, 'Primera', 750, 14000, 15000, 14600, '$/caja 13 kilos', 1123, 13),
        @(28, 45049, 'Primera', 300, 13000, 14000, 13500, '$/caja 13 kilos', 1038, 13)
)

foreach ($r in $rowsData) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 4).Value  = $r[1]   # D - Fecha
    $ws.Cells.Item($rowNum, 9).Value  = $r[2]   # I - Calidad
    $ws.Cells.Item($rowNum, 10).Value = $r[3]   # J - Volumen
    $ws.Cells.Item($rowNum, 11).Value = $r[4]   # K - Precio minimo
    $ws.Cells.Item($rowNum, 12).Value = $r[5]   # L - Precio maximo
    $ws.Cells.Item($rowNum, 13).Value = $r[6]   # M - Precio promedio ponderado
    $ws.Cells.Item($rowNum, 14).Value = $r[7]   # N - Unidad de comercializacion
    $ws.Cells.Item($rowNum, 16).Value = $r[8]   # P - Precio $/Kg
    $ws.Cells.Item($rowNum, 17).Value = $r[9]   # Q - Kg o Unidades
}
